# "Use <formatversion> as suffix for table headers"
#
# The sheet holds a diff table comparing two format-versions:
#   columns A..J   -> header names ending in "_old"   => rename suffix to "_FV2210"
#   column  K      -> "diff"                          => unchanged
#   columns L..U   -> header names ending in "_new"   => rename suffix to "_FV2304"
#
# Afterwards the used range A1:U56 becomes a proper Excel Table ("Table1")
# and the header row (row 1) is frozen.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 21   # column U
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = [string]$cell.Text

    if ($header.EndsWith("_old")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2210"
    }
    elseif ($header.EndsWith("_new")) {
        $cell.Value = $header.Substring(0, $header.Length - 4) + "_FV2304"
    }
}

# --- Turn the used range into an Excel Table ---------------------------
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U56"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- Freeze the header row ----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
